$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 123
$ws1.Range("F4").Value = 159
$ws1.Range("F5").Value = 3161
$ws1.Range("F6").Value = 317
$ws1.Range("F7").Value = 8
$ws1.Range("F8").Value = 414

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 123
$ws4.Range("F4").Value = 159
$ws4.Range("F5").Value = 3161
$ws4.Range("F6").Value = 317
$ws4.Range("F9").Value = 8
$ws4.Range("F10").Value = 414
